# Auto-committed on 2022/09/16 週五 17:39:19.43
# Inserts 3 new field rows ("ReconCode"/"對帳類別", "TitaTlrNo"/"經辦",
# "TitaTxtNo"/"交易序號") into the DBD field list of BankRmtf.xlsx, right
# before the existing "RemintBank" row, renumbering the SEQ column for the
# rows pushed down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# ---------------------------------------------------------------------
# 1. Make room: insert 3 blank rows at 25..27 (shifts old rows 25-28 -> 28-31)
# ---------------------------------------------------------------------
$ws.Rows("25:27").Insert(-4121, 0)   # xlShiftDown, xlFormatFromLeftOrAbove

# ---------------------------------------------------------------------
# 2. Copy the cell formatting that the new rows should use from existing
#    rows/cells with matching look & feel, then fill in content.
# ---------------------------------------------------------------------

# --- Row 26: SEQ 17 / TitaTlrNo / 經辦 / VARCHAR2 / 6 -------------------
$ws.Range("A9").Copy()
$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B9").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("A26").Value = 17
$ws.Range("B26").Value = "TitaTlrNo"
$ws.Range("C26").Value = "經辦"
$ws.Range("D26").Value = "VARCHAR2"
$ws.Range("E26").Value = 6

# --- Row 27: SEQ 18 / TitaTxtNo / 交易序號 / VARCHAR2 / 8 ---------------
$ws.Range("A9").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G27").PasteSpecial(-4122)

$ws.Range("A27").Value = 18
$ws.Range("B27").Value = "TitaTxtNo"
$ws.Range("C27").Value = "交易序號"
$ws.Range("D27").Value = "VARCHAR2"
$ws.Range("E27").Value = 8

# --- Row 25: SEQ 16 / ReconCode / 對帳類別 / VARCHAR2 / 3 + long note ---
$ws.Range("A9").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("G24").Copy()
$ws.Range("G25").PasteSpecial(-4122)

$ws.Range("A25").Value = 16
$ws.Range("B25").Value = "ReconCode"
$ws.Range("C25").Value = "對帳類別"
$ws.Range("D25").Value = "VARCHAR2"
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = "轉AS400匯款轉帳檔的存摺代號(舊資料有01,02,03,13亦新增代碼但目前已不使用)`nCdCode.ReconCode`nP01:銀行存款－郵局`nC01:暫收款－非核心運用`nA1~A7:  (帳務:P03:銀行存款－新光)`nTEM:員工扣薪15/非15`nTCK:支票"

$ws.Rows("25").RowHeight = 129.6

# ---------------------------------------------------------------------
# 3. Renumber the SEQ column for the rows that were pushed down
#    (old 25-28, now 28-31): 16,17,18,19 -> 19,20,21,22
# ---------------------------------------------------------------------
$ws.Range("A28").Value = 19
$ws.Range("A29").Value = 20
$ws.Range("A30").Value = 21
$ws.Range("A31").Value = 22

# ---------------------------------------------------------------------
# 4. Restore the on-screen selection/scroll position to roughly match
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("H25").Select()

$wb.Save()
